# Route to get LDJSON, added more fields
$wb = $excel.ActiveWorkbook

# Work on the "routes" worksheet (sheet2 in the package)
$ws = $wb.Worksheets.Item("routes")
$ws3 = $wb.Worksheets.Item("localizedRegex")

# Match the style used by the existing data rows (E5/G5 stay blank but keep
# the row's formatting)
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Add the new row (row 5) of data to the routes table.
# Values are entered in the same order the author typed them so the shared
# strings table grows in the matching sequence (Params, Function,
# Description, then Name).
$ws.Cells.Item(5, 1).Value = "recipe"
$ws.Cells.Item(5, 4).Value = "url*"
$ws.Cells.Item(5, 3).Value = "getLDJSONfromURL"
$ws.Cells.Item(5, 6).Value = "Gets LD JSON from a url"
$ws.Cells.Item(5, 2).Value = "Gets LD JSON"

# Resize the table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Table4")
$table.Resize($ws.Range("A2:G5"))

# Widen column C to fit the new, longer "getLDJSONfromURL" entry
$ws.Columns.Item(3).ColumnWidth = 17.6

# The routes sheet becomes the active tab/sheet, with B6 selected
$ws.Activate()
$ws.Range("B6").Select()

# The previously active sheet (localizedRegex) keeps its own (new) selection
$ws3.Range("A12").Select()

# Re-activate routes so it ends up as the active tab (matches activeTab=1)
$ws.Activate()
